# TrialsSetup update — 2026-02-24 16:00
# REMASTER (CLOU): Days remaining 16 -> 15
# COLO-PREVENT:    Progress        0 -> 12.5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 = REMASTER (CLOU), column B = "Days remaining"
$ws.Range("B8").Value = 15

# Row 9 = COLO-PREVENT, column C = "Progress"
$ws.Range("C9").Value = 12.5
